$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style to the built-in table style used after the edit.
# ---------------------------------------------------------------------------
$newTableStyleId = "{45B174DE-356D-4AD1-B9D7-024270EB1710}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's theme palette from "Integral" (Red Violet) to
#    the default "Office Theme" colours.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hyperlink
    0x954F72   # 12 followed hyperlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $rVal = ($officeColors[$i - 1] -band 0xFF0000) -shr 16
    $gVal = ($officeColors[$i - 1] -band 0x00FF00) -shr 8
    $bVal = ($officeColors[$i - 1] -band 0x0000FF)
    $comRgb = $rVal -bor ($gVal -shl 8) -bor ($bVal -shl 16)
    $tcs.Colors($i).RGB = $comRgb
}
